# Auto-generated script to update cached market-price figures in each class sheet
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2392.0715
$ws.Range("I40").Value = 2264.75
$ws.Range("J40").Value = 2561.8333
$ws.Range("K40").Value = 2264.75
$ws.Range("L40").Value = 2561.8333
$ws.Range("M40").Value = -2089.75
$ws.Range("N40").Value = -2911.8333
$ws.Range("H88").Value = 7517
$ws.Range("J88").Value = 8020.4
$ws.Range("L88").Value = 8020.4
$ws.Range("N88").Value = -8832.4
$ws.Range("H91").Value = 7517
$ws.Range("J91").Value = 8020.4
$ws.Range("L91").Value = 8020.4
$ws.Range("N91").Value = -10828.4
$ws.Range("H100").Value = 3527.8572
$ws.Range("I100").Value = 3539.2
$ws.Range("J100").Value = 3499.5
$ws.Range("K100").Value = 3539.2
$ws.Range("L100").Value = 3499.5
$ws.Range("M100").Value = -2998.2
$ws.Range("N100").Value = -4581.5
$ws.Range("H106").Value = 32000
$ws.Range("I106").Value = 34000
$ws.Range("K106").Value = 34000
$ws.Range("M106").Value = -33369
$ws.Range("H112").Value = 3819.0625
$ws.Range("J112").Value = 3039.1538
$ws.Range("L112").Value = 9117.4614
$ws.Range("N112").Value = -11333.4614
$ws.Range("H127").Value = 10192.667
$ws.Range("I127").Value = 10192.667
$ws.Range("K127").Value = 30578.001
$ws.Range("M127").Value = -25618.001
$ws.Range("H131").Value = 865
$ws.Range("I131").Value = 865
$ws.Range("K131").Value = 2595
$ws.Range("M131").Value = 2445
$ws.Range("H132").Value = 7244.5
$ws.Range("I132").Value = 7693.8
$ws.Range("K132").Value = 23081.4
$ws.Range("M132").Value = -20551.4
$ws.Range("H135").Value = 1965
$ws.Range("I135").Value = 433
$ws.Range("K135").Value = 3897
$ws.Range("M135").Value = -1362
$ws.Range("H137").Value = 2137.6
$ws.Range("I137").Value = 1962.6666
$ws.Range("K137").Value = 5887.9998
$ws.Range("M137").Value = -3337.9998
$ws.Range("H138").Value = 1486.762
$ws.Range("I138").Value = 1220.2354
$ws.Range("J138").Value = 2619.5
$ws.Range("K138").Value = 3660.7062
$ws.Range("L138").Value = 7858.5
$ws.Range("M138").Value = 1479.2938
$ws.Range("N138").Value = -18138.5
$ws.Range("H141").Value = 18031.334
$ws.Range("I141").Value = 17999
$ws.Range("K141").Value = 53997
$ws.Range("M141").Value = -48817

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 15000
$ws.Range("J42").Value = 15000
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = -15972
$ws.Range("H122").Value = 1049.5714
$ws.Range("I122").Value = 1049.5714
$ws.Range("K122").Value = 3148.7142
$ws.Range("M122").Value = -698.7142000000003
$ws.Range("H132").Value = 1745.7333
$ws.Range("I132").Value = 1730.8636
$ws.Range("J132").Value = 1786.625
$ws.Range("K132").Value = 5192.5908
$ws.Range("L132").Value = 5359.875
$ws.Range("M132").Value = -2662.5908
$ws.Range("N132").Value = -10419.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 303.81818
$ws.Range("I22").Value = 326.55554
$ws.Range("J22").Value = 201.5
$ws.Range("K22").Value = 326.55554
$ws.Range("L22").Value = 201.5
$ws.Range("M22").Value = -153.55554
$ws.Range("N22").Value = -547.5
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H94").Value = 433.33334
$ws.Range("I94").Value = 433.33334
$ws.Range("K94").Value = 433.33334
$ws.Range("M94").Value = 17.66665999999998
$ws.Range("H99").Value = 26145964
$ws.Range("I99").Value = 8548616
$ws.Range("J99").Value = 83337350
$ws.Range("K99").Value = 8548616
$ws.Range("L99").Value = 83337350
$ws.Range("M99").Value = -8547118
$ws.Range("N99").Value = -83340346
$ws.Range("H134").Value = 3529.9092
$ws.Range("I134").Value = 3870.2222
$ws.Range("J134").Value = 1998.5
$ws.Range("K134").Value = 11610.6666
$ws.Range("L134").Value = 5995.5
$ws.Range("M134").Value = -9075.6666
$ws.Range("N134").Value = -11065.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 604.75
$ws.Range("J7").Value = 413
$ws.Range("L7").Value = 413
$ws.Range("N7").Value = -639
$ws.Range("H22").Value = 882.9286
$ws.Range("I22").Value = 931.6667
$ws.Range("J22").Value = 846.375
$ws.Range("K22").Value = 931.6667
$ws.Range("L22").Value = 846.375
$ws.Range("M22").Value = -581.6667
$ws.Range("N22").Value = -1546.375
$ws.Range("H31").Value = 4939.4
$ws.Range("I31").Value = 4174.25
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 4174.25
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = -3879.25
$ws.Range("N31").Value = -8590
$ws.Range("H34").Value = 4939.4
$ws.Range("I34").Value = 4174.25
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 4174.25
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = -3972.25
$ws.Range("N34").Value = -8404
$ws.Range("H37").Value = 24494.5
$ws.Range("I37").Value = 24494.5
$ws.Range("K37").Value = 24494.5
$ws.Range("M37").Value = -24387.5
$ws.Range("H58").Value = 3092.5334
$ws.Range("I58").Value = 2613.6365
$ws.Range("K58").Value = 2613.6365
$ws.Range("M58").Value = -2410.6365
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H99").Value = 3708.6
$ws.Range("I99").Value = 4181.6665
$ws.Range("J99").Value = 2999
$ws.Range("K99").Value = 4181.6665
$ws.Range("L99").Value = 2999
$ws.Range("M99").Value = -2683.6665
$ws.Range("N99").Value = -5995
$ws.Range("H126").Value = 3708.6
$ws.Range("I126").Value = 4181.6665
$ws.Range("J126").Value = 2999
$ws.Range("K126").Value = 12544.9995
$ws.Range("L126").Value = 8997
$ws.Range("M126").Value = -10074.9995
$ws.Range("N126").Value = -13937
$ws.Range("H132").Value = 2283.4707
$ws.Range("I132").Value = 1943.6666
$ws.Range("K132").Value = 5830.9998
$ws.Range("M132").Value = -3300.9998
$ws.Range("H136").Value = 3092.5334
$ws.Range("I136").Value = 2613.6365
$ws.Range("K136").Value = 7840.9095
$ws.Range("M136").Value = -5290.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 73.8
$ws.Range("I2").Value = 69.666664
$ws.Range("K2").Value = 417.999984
$ws.Range("M2").Value = -304.999984
$ws.Range("H113").Value = 1239.7273
$ws.Range("J113").Value = 1348.2222
$ws.Range("L113").Value = 4044.6666
$ws.Range("N113").Value = -8384.6666
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H136").Value = 2549.75
$ws.Range("I136").Value = 2549.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7649.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2549.25
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 5533.25
$ws.Range("I137").Value = 2600
$ws.Range("K137").Value = 7800
$ws.Range("M137").Value = -2700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 46633
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 46633
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 46633
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -48005
$ws.Range("H66").Value = 46633
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 46633
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 139899
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -146763
$ws.Range("H97").Value = 334.27274
$ws.Range("I97").Value = 386.66666
$ws.Range("K97").Value = 386.66666
$ws.Range("M97").Value = 109.33334
$ws.Range("H126").Value = 4444
$ws.Range("I126").Value = 4444
$ws.Range("K126").Value = 13332
$ws.Range("M126").Value = -10862
$ws.Range("H132").Value = 3527.1428
$ws.Range("I132").Value = 3365.0833
$ws.Range("K132").Value = 10095.2499
$ws.Range("M132").Value = -7565.249899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1900
$ws.Range("I100").Value = 1900
$ws.Range("K100").Value = 1900
$ws.Range("M100").Value = -1359
$ws.Range("H122").Value = 21249.75
$ws.Range("J122").Value = 20000
$ws.Range("L122").Value = 60000
$ws.Range("N122").Value = -64900
$ws.Range("H132").Value = 2876.625
$ws.Range("I132").Value = 2834.3333
$ws.Range("J132").Value = 3003.5
$ws.Range("K132").Value = 8502.999899999999
$ws.Range("L132").Value = 9010.5
$ws.Range("M132").Value = -5972.999899999999
$ws.Range("N132").Value = -14070.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 98996
$ws.Range("J121").Value = 98996
$ws.Range("L121").Value = 98996
$ws.Range("N121").Value = -102490
$ws.Range("H126").Value = 2143
$ws.Range("I126").Value = 2166.8333
$ws.Range("K126").Value = 6500.499899999999
$ws.Range("M126").Value = -4030.499899999999
$ws.Range("H136").Value = 6634
$ws.Range("I136").Value = 7501
$ws.Range("J136").Value = 4900
$ws.Range("K136").Value = 22503
$ws.Range("L136").Value = 14700
$ws.Range("M136").Value = -19953
$ws.Range("N136").Value = -19800
